$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (shifting C..W to D..X).
# Excel's column-insert inherits the formatting of the column to the left,
# so the new header cell automatically picks up the bold/fill/border style
# used by the rest of row 1.
$ws.Columns("C").Insert()

# Set the header for the newly inserted column C1
$ws.Range("C1").Value = "Legal Entity"

# New column inherits the width of column B (same as Excel's column insert behavior)
$ws.Columns("C").ColumnWidth = $ws.Columns("B").ColumnWidth

# Update selection to match the target state
$ws.Range("D13").Select()
